# Apply updated HIGH/LOW/CLOSE/LTP/VOL/9:25 CLOSE values to Sheet1
# (new algo path added to MWCL creation updation.py refreshed these figures)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 2).Value = 1416.95
$ws.Cells.Item(2, 3).Value = 1381.25
$ws.Cells.Item(2, 4).Value = 3004.7
$ws.Cells.Item(2, 5).Value = 1410
$ws.Cells.Item(2, 6).Value = 33
$ws.Cells.Item(2, 7).Value = 1386.35

# Row 3
$ws.Cells.Item(3, 2).Value = 1488.1
$ws.Cells.Item(3, 3).Value = 1455.35
$ws.Cells.Item(3, 4).Value = 1478.4
$ws.Cells.Item(3, 5).Value = 1479.25
$ws.Cells.Item(3, 6).Value = 16
$ws.Cells.Item(3, 7).Value = 1459.65

# Row 4
$ws.Cells.Item(4, 2).Value = 52420
$ws.Cells.Item(4, 3).Value = 51222.35
$ws.Cells.Item(4, 4).Value = 52361.4
$ws.Cells.Item(4, 5).Value = 52310.4
$ws.Cells.Item(4, 6).Value = 25
$ws.Cells.Item(4, 7).Value = 51250.05

# Row 5
$ws.Cells.Item(5, 2).Value = 105.5
$ws.Cells.Item(5, 3).Value = 101.01
$ws.Cells.Item(5, 4).Value = 105
$ws.Cells.Item(5, 5).Value = 105.12
$ws.Cells.Item(5, 6).Value = 589
$ws.Cells.Item(5, 7).Value = 101.09

# Row 6
$ws.Cells.Item(6, 2).Value = 882.5
$ws.Cells.Item(6, 3).Value = 851.1
$ws.Cells.Item(6, 4).Value = 877
$ws.Cells.Item(6, 5).Value = 878.4
$ws.Cells.Item(6, 6).Value = 89
$ws.Cells.Item(6, 7).Value = 853.05

# Row 7
$ws.Cells.Item(7, 2).Value = 759.5
$ws.Cells.Item(7, 3).Value = 727.45
$ws.Cells.Item(7, 4).Value = 756.65
$ws.Cells.Item(7, 5).Value = 756.9
$ws.Cells.Item(7, 6).Value = 143
$ws.Cells.Item(7, 7).Value = 727.5

# Row 8
$ws.Cells.Item(8, 2).Value = 1269.3
$ws.Cells.Item(8, 3).Value = 1229.55
$ws.Cells.Item(8, 4).Value = 1268.95
$ws.Cells.Item(8, 5).Value = 1266.85
$ws.Cells.Item(8, 6).Value = 287
$ws.Cells.Item(8, 7).Value = 1231.95

# Row 9
$ws.Cells.Item(9, 2).Value = 975
$ws.Cells.Item(9, 3).Value = 925.4
$ws.Cells.Item(9, 4).Value = 968.75
$ws.Cells.Item(9, 5).Value = 968.35
$ws.Cells.Item(9, 6).Value = 53
$ws.Cells.Item(9, 7).Value = 925.9

# Row 10
$ws.Cells.Item(10, 2).Value = 24996.75
$ws.Cells.Item(10, 3).Value = 24640.2
$ws.Cells.Item(10, 4).Value = 24965.55
$ws.Cells.Item(10, 5).Value = 24949.15
$ws.Cells.Item(10, 6).Value = 61
$ws.Cells.Item(10, 7).Value = 24656.5

# Row 11
$ws.Cells.Item(11, 2).Value = 2744.95
$ws.Cells.Item(11, 3).Value = 2692
$ws.Cells.Item(11, 4).Value = 2727
$ws.Cells.Item(11, 5).Value = 2728.6
$ws.Cells.Item(11, 6).Value = 120
$ws.Cells.Item(11, 7).Value = 2693.2

# Row 12
$ws.Cells.Item(12, 2).Value = 826.7
$ws.Cells.Item(12, 3).Value = 805.05
$ws.Cells.Item(12, 4).Value = 823.4
$ws.Cells.Item(12, 5).Value = 824.3
$ws.Cells.Item(12, 6).Value = 309
$ws.Cells.Item(12, 7).Value = 806.7

# Row 13
$ws.Cells.Item(13, 2).Value = 1102.4
$ws.Cells.Item(13, 3).Value = 1076.7
$ws.Cells.Item(13, 4).Value = 1097.45
$ws.Cells.Item(13, 5).Value = 1098.3
$ws.Cells.Item(13, 6).Value = 28
$ws.Cells.Item(13, 7).Value = 1080.6

# Row 14
$ws.Cells.Item(14, 2).Value = 920.5
$ws.Cells.Item(14, 3).Value = 891.6
$ws.Cells.Item(14, 4).Value = 913.8
$ws.Cells.Item(14, 5).Value = 913.55
$ws.Cells.Item(14, 6).Value = 163
$ws.Cells.Item(14, 7).Value = 892.65

# Row 15
$ws.Cells.Item(15, 2).Value = 156.9
$ws.Cells.Item(15, 3).Value = 151.62
$ws.Cells.Item(15, 4).Value = 156.01
$ws.Cells.Item(15, 5).Value = 156.09
$ws.Cells.Item(15, 6).Value = 642
$ws.Cells.Item(15, 7).Value = 151.7

# Row 16
$ws.Cells.Item(16, 2).Value = 4152.95
$ws.Cells.Item(16, 3).Value = 4111.65
$ws.Cells.Item(16, 4).Value = 4140
$ws.Cells.Item(16, 5).Value = 4136.3
$ws.Cells.Item(16, 6).Value = 32
$ws.Cells.Item(16, 7).Value = 4129.05

# Row 17
$ws.Cells.Item(17, 2).Value = 3429.65
$ws.Cells.Item(17, 3).Value = 3301.75
$ws.Cells.Item(17, 4).Value = 3408.7
$ws.Cells.Item(17, 5).Value = 3397.5
$ws.Cells.Item(17, 6).Value = 36
$ws.Cells.Item(17, 7).Value = 3340.55
